# "delete iteration option, added error handling for user inputs - Reem"
# Populate the user-story table (columns A:name, B:size, C:iteration) with
# the new rows. Cell-by-cell order matches the order new shared strings were
# introduced in the target workbook: column A first for the "seed" rows,
# then column C for those same rows, then column A for the rows that were
# inserted afterwards, and finally the remaining numeric/text fill-ins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: column A for the originally-entered rows (2,3,5,7) ---
$ws.Range("A2").Value = "US11"
$ws.Range("A3").Value = "US12"
$ws.Range("A5").Value = "US22"
$ws.Range("A7").Value = "US31"

# --- Step 2: column C (iteration) for those same rows ---
$ws.Range("C2").Value = "Itr 1"
$ws.Range("C3").Value = "Itr 1"
$ws.Range("C5").Value = "Itr 2"
$ws.Range("C7").Value = "Itr 3"

# --- Step 3: column A for the rows added later (4,6,8) ---
$ws.Range("A4").Value = "US21"
$ws.Range("A6").Value = "US23"
$ws.Range("A8").Value = "US32"

# --- Step 4: remaining column C entries for the newer rows ---
$ws.Range("C4").Value = "Itr 2"
$ws.Range("C6").Value = "Itr 2"
$ws.Range("C8").Value = "Itr 3"

# --- Step 5: fill in the "size" column (B) for every new row ---
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 8

# Final selection left on C9, matching the saved workbook state.
$ws.Range("C9").Select()
